$wb = $excel.ActiveWorkbook

# --- Reorder worksheet tabs: move "Add Exchanges" to sit before "Delete Exchanges" ---
$wsAdd = $wb.Worksheets.Item("Add Exchanges")
$wsDelete = $wb.Worksheets.Item("Delete Exchanges")
$wsAdd.Move($wsDelete)

# --- Update the two "Add Exchanges" rows whose Activity Database / Activity
#     Reference Code values changed (new database/code from new act -> newdb/uuid) ---
# Re-fetch the worksheet reference post-move so it is not stale.
$wsAdd = $wb.Worksheets.Item("Add Exchanges")
$wsAdd.Range("A4").Value = "newdb"
$wsAdd.Range("C4").Value = "uuid"
$wsAdd.Range("A5").Value = "newdb"
$wsAdd.Range("C5").Value = "uuid"

# --- Update selection / active cell on "Add Exchanges" and make it the active sheet ---
$wsAdd.Activate()
$wsAdd.Range("C8").Select()

# --- Update selection / active cell on "Create Activities" ---
$wsCreate = $wb.Worksheets.Item("Create Activities")
$wsCreate.Range("J1").Select()

# --- Leave "Add Exchanges" as the active / displayed sheet, matching the target file ---
$wsAdd.Activate()
